{"js": "const body = context.document.body;\n\n// 1) \"June\" -> \"November\" in the Date paragraph.\nconst found = body.search(\"June\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"November\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Collapse the seven figure paragraphs (FirstParagraph + six BodyText)\n//    into a single paragraph (keeping the FirstParagraph style) whose text\n//    reads \"results\", removing the inline pictures.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nconst figureParas = paragraphs.items.filter((p) => p.style === \"First Paragraph\" || p.style === \"Body Text\");\n\n// Keep the first figure paragraph (FirstParagraph style) to host the new\n// text, delete the rest.\nconst keepPara = figureParas[0];\nfor (let i = 1; i < figureParas.length; i++) {\n  figureParas[i].delete();\n}\nawait context.sync();\n\n// Remove the picture from the kept paragraph and add the \"results\" text.\nconst keepPics = keepPara.inlinePictures;\nkeepPics.load(\"items\");\nawait context.sync();\nfor (const pic of keepPics.items) {\n  pic.delete();\n}\nkeepPara.insertText(\"results\", \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"June\" -> \"November\" in the Date paragraph.\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Replacement.ClearFormatting()\n$range.Find.Execute(\"June\", $false, $false, $false, $false, $false, $true, 1, $false, \"November\", 2) | Out-Null\n\n# 2) Collapse the seven figure paragraphs (FirstParagraph + six BodyText)\n#    into a single paragraph (keeping the FirstParagraph style) whose text\n#    reads \"results\", removing the inline pictures.\n$total = $d.Paragraphs.Count\nfor ($i = $total; $i -ge 4; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\nfor ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {\n    $d.InlineShapes.Item($i).Delete()\n}\n\n$keepPara = $d.Paragraphs.Item(3)\n$keepPara.Range.Text = \"results\"\n"}
